$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Unmerge the old group-header cells (H1:L1, M1:P1, Q1:S1) first so each
#    individual cell underneath can hold its own value again.
$ws.Range("H1:L1").UnMerge()
$ws.Range("M1:P1").UnMerge()
$ws.Range("Q1:S1").UnMerge()

# 2. Flatten the two-row MultiIndex header into row 1: replace the old
#    "Unnamed: N_level_0" / group-title placeholders with the real column
#    labels (mirroring row 2), reordering Tkl% / Lost / Blocks / Sh, and add
#    the extra identifier columns (Player ID / 90s / Cha).
$ws.Range("A1").Value = "Player ID"
$ws.Range("B1").Value = "Player"
$ws.Range("C1").Value = "#"
$ws.Range("D1").Value = "Nation"
$ws.Range("E1").Value = "Pos"
$ws.Range("F1").Value = "Age"
$ws.Range("G1").Value = "90s"
$ws.Range("H1").Value = "Tkl"
$ws.Range("I1").Value = "TklW"
$ws.Range("J1").Value = "Def 3rd"
$ws.Range("K1").Value = "Mid 3rd"
$ws.Range("L1").Value = "Att 3rd"
$ws.Range("M1").Value = "Cha"
$ws.Range("N1").Value = "Att"
$ws.Range("O1").Value = "Tkl%"
$ws.Range("P1").Value = "Lost"
$ws.Range("Q1").Value = "Blocks"
$ws.Range("R1").Value = "Sh"
$ws.Range("S1").Value = "Pass"
$ws.Range("T1").Value = "Int"
$ws.Range("U1").Value = "Tkl+Int"
$ws.Range("V1").Value = "Clr"
$ws.Range("W1").Value = "Err"

# 3. The original column-label row (row 2) is now redundant - keep it for
#    reference but hide it, along with the blank spacer row 3 beneath it,
#    and hide the totals/summary row 20 now that it's just a footer.
$ws.Rows.Item(2).Hidden = $true
$ws.Rows.Item(3).Hidden = $true
$ws.Rows.Item(20).Hidden = $true

# 4. Fill in explicit 0s for the Tkl% cells that were left blank for players
#    with no tackle attempts.
$ws.Range("O4").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("O6").Value = 0
$ws.Range("O9").Value = 0
$ws.Range("O16").Value = 0
$ws.Range("O17").Value = 0
$ws.Range("O18").Value = 0
$ws.Range("O19").Value = 0

# 5. Leave the selection where the cleanup left off.
$ws.Range("O21").Select() | Out-Null
